# The "provincia" column (G) was previously modeled as an sdmx dimension
# (refArea / dim / URI-Provincia). The data was reprocessed with the new
# curated dimensions so it is now modeled like the other iaest measures
# (iaest-measure:provincia / medida / xsd:int). This also drops the now
# unused "URI-Provincia" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "iaest-measure:provincia"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
